# Split the "Normal Distributions" topic (which spanned two identical rows)
# into two distinct topics: an "Introduction" and a "Calculations" session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "NormalDist1"
$ws.Range("E9").Value = "NormalDist2"

$ws.Range("D8").Value = "Normal Distributions - Introduction"
$ws.Range("D9").Value = "Normal Distributions - Calculations"

# Column D needs to widen to fit the longer new text (auto-fit like Excel would do).
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668

# Reflect the cell the editor ended up on after making the change.
$ws.Range("G16").Select() | Out-Null
